$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,2,5,1),
    @(3,2,4,0),
    @(5,0,6,2),
    @(3,2,4,1),
    @(7,3,4,0),
    @(4,2,4,1),
    @(4,2,6,0),
    @(4,2,3,1),
    @(6,2,7,1),
    @(2,2,3,0),
    @(4,1,4,2),
    @(3,0,4,3),
    @(5,2,5,0),
    @(6,3,6,0),
    @(2,3,3,0),
    @(3,1,4,2),
    @(4,2,4,1),
    @(6,3,5,0),
    @(5,2,5,0),
    @(4,2,5,0),
    @(4,2,4,0),
    @(4,2,3,1),
    @(6,2,4,1),
    @(2,2,3,1),
    @(3,1,3,2),
    @(5,2,5,0),
    @(6,0,5,2),
    @(6,0,7,2),
    @(3,3,3,0),
    @(4,2,4,1),
    @(6,2,5,1),
    @(3,1,4,2),
    @(5,0,6,3),
    @(5,2,2,1),
    @(3,3,2,0),
    @(3,2,3,1),
    @(6,2,4,0),
    @(7,2,7,0),
    @(3,1,5,2),
    @(3,3,3,0),
    @(5,1,5,2),
    @(4,1,4,2),
    @(6,0,6,2),
    @(4,2,5,0)
)

$startRow = 3256
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
}

$ws.Range("A3300").Select() | Out-Null
